$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# The "Greece" sheet is the template used for every per-market tab in this
# workbook (header block in rows 1-7, then an "Attached Functionality" list
# below). Three new market tabs are added the same way: copy the template,
# rename it, and fill in the market-specific cells.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("Greece")

# 1) Netherlands - placed right after Greece
$template.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$nl = $wb.Worksheets.Item($wb.Worksheets.Count)
$nl.Name = "Netherlands"

# 2) Austria - placed right after Netherlands
$template.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$at = $wb.Worksheets.Item($wb.Worksheets.Count)
$at.Name = "Austria"

# 3) Denmark - placed right after Austria
$template.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$dk = $wb.Worksheets.Item($wb.Worksheets.Count)
$dk.Name = "Denmark"

# ---------------------------------------------------------------------------
# Fill in the "Market" name cells first (matches shared-string insert order)
# ---------------------------------------------------------------------------
$nl.Range("B2").Value = "Netherlands Market"
$at.Range("B2").Value = "Austria Market"

# ---------------------------------------------------------------------------
# Austria has extra product rows: it uses a "Multichannel Transmission Unit"
# in addition to the plain "Transmission Unit" / "Transmission Unit and
# Keysafe" combos, so 5 extra rows are inserted before the trailing
# Black Box / Wg / Attached Functionality rows.
# ---------------------------------------------------------------------------
$at.Range("A10:A14").EntireRow.Insert()
$at.Range("A9").Copy()
$at.Range("A10:A14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$at.Application.CutCopyMode = $false

$at.Range("A10").Value = "Multichannel Transmission Unit"
$at.Range("A11").Value = "Transmission Unit"
$at.Range("A12").Value = "Transmission Unit and Keysafe"
$at.Range("A13").Value = "Transmission Unit and Keysafe"
$at.Range("A14").Value = "Multichannel Transmission Unit"

# ---------------------------------------------------------------------------
# Fill in the "Attached Functionality" reference (NGC/T) cells
# ---------------------------------------------------------------------------
$nl.Range("B4").Value = "NGC-3144/T2183"
$at.Range("B4").Value = "NGC-3817/T2279"

$dk.Range("B2").Value = "Denmark Market"
$dk.Range("B4").Value = "NGC-2913/T2279"

# ---------------------------------------------------------------------------
# Re-size the columns on the new tabs (wider text -> wider columns)
# ---------------------------------------------------------------------------
$nl.Columns.Item(1).ColumnWidth = 42
$nl.Columns.Item(2).ColumnWidth = 38
$nl.Columns.Item(3).ColumnWidth = 12.333333333333332
$nl.Columns.Item(4).ColumnWidth = 22

$at.Columns.Item(1).ColumnWidth = 42.5
$at.Columns.Item(2).ColumnWidth = 25.333333333333336
$at.Columns.Item(3).ColumnWidth = 13.666666666666666
$at.Columns.Item(4).ColumnWidth = 16.666666666666668

$dk.Columns.Item(1).ColumnWidth = 42.5
$dk.Columns.Item(2).ColumnWidth = 25.333333333333336
$dk.Columns.Item(3).ColumnWidth = 13.666666666666666
$dk.Columns.Item(4).ColumnWidth = 16.666666666666668

# ---------------------------------------------------------------------------
# Restore each sheet's own selection state
# ---------------------------------------------------------------------------
$nl.Range("A8:A12").Select()
$at.Cells.Select()
$dk.Range("A15").Select()

# Denmark ends up as the active/selected tab
$dk.Activate()
